$wb = $excel.ActiveWorkbook

# NOTE: in the source workbook, every numeric-looking value in the
# Restricciones_del_follower / Punto_modificado / Vector_bf / Vector_BF
# sheets is actually stored as TEXT (shared string), not as a real
# number (the file was produced by a non-Excel writer). Assigning a
# numeric-looking string straight to .Value auto-coerces it to a
# number, so we force the cell format to Text ("@") first for those
# cells to keep them as text, matching the original data model. The
# Vector_Alpha sheet, by contrast, genuinely stores real numbers there
# already, so those are left as plain numeric assignments.

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
}

# --- Sheet: Restricciones_del_follower ---
$ws = $wb.Worksheets.Item("Restricciones_del_follower")

# Row 2
$ws.Range("A2").Value = "-0.9446013643891683y_1 + 0.011004345402239183y_2"
Set-TextValue $ws.Range("B2") "-4.113656050722677"
$ws.Range("C2").Value = "J_0_L0_v"
Set-TextValue $ws.Range("D2") "0.0866877650392671"
Set-TextValue $ws.Range("E2") "-0.2704755559527997"
Set-TextValue $ws.Range("F2") "-0.4724490455978724"

# Row 3
$ws.Range("A3").Value = "-4 + 1.4258984934749632y_1 + 0.08460017249189158y_2"
Set-TextValue $ws.Range("B3") "2.4513323248895706"
$ws.Range("C3").Value = "J_0_L0_v"
Set-TextValue $ws.Range("D3") "0.9648587319705634"
Set-TextValue $ws.Range("E3") "0"
Set-TextValue $ws.Range("F3") "0.7232519754967422"

# Row 4
$ws.Range("A4").Value = "-16 - 2x + 0.11759800913107077y_1 + 3.8247202989012488y_2"
Set-TextValue $ws.Range("B4") "-18.10373805947259"
$ws.Range("C4").Value = "J_0_LP_v"
Set-TextValue $ws.Range("D4") "0.9761226555169311"
Set-TextValue $ws.Range("E4") "0.4765093561492265"
Set-TextValue $ws.Range("F4") "0"

# Row 5
$ws.Range("A5").Value = "-48 + 8x + 1.9304529951548228y_1 + 0.18482451826358154y_2"
Set-TextValue $ws.Range("B5") "7.908693258692484"
$ws.Range("C5").Value = "J_Ne_L0_v"
Set-TextValue $ws.Range("D5") "0.8143958706897286"
Set-TextValue $ws.Range("E5") "0.9523557914603366"
Set-TextValue $ws.Range("F5") "0.9861742737355413"

# Row 6
$ws.Range("A6").Value = "12 - 2x + 0.8716114642162638y_1 + 0.5704148498394974y_2"
Set-TextValue $ws.Range("B6") "5.430374185532117"
$ws.Range("C6").Value = "J_Ne_L0_v"
Set-TextValue $ws.Range("D6") "0.23927405565041526"
Set-TextValue $ws.Range("E6") "0"
Set-TextValue $ws.Range("F6") "0.49375642524538843"

# --- Sheet: Punto_modificado ---
$ws4 = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $ws4.Range("A2") "5.875840352759835"
Set-TextValue $ws4.Range("B2") "4.382729079133727"
Set-TextValue $ws4.Range("C2") "2.387767396848251"

# --- Sheet: Vector_bf (index 5, since name lookup is case-insensitive
#     and would collide with Vector_BF) ---
$ws5 = $wb.Worksheets.Item(5)
Set-TextValue $ws5.Range("A2") "-2.189402269966704"
Set-TextValue $ws5.Range("A3") "-4.1029830910370535"

# --- Sheet: Vector_BF (index 6) ---
$ws6 = $wb.Worksheets.Item(6)
Set-TextValue $ws6.Range("A2") "-5.665827619384239"
Set-TextValue $ws6.Range("A3") "0.8499937791199406"
Set-TextValue $ws6.Range("A4") "-3.9955573010119325"

# --- Sheet: Vector_Alpha (values here are real numbers already) ---
$ws7 = $wb.Worksheets.Item("Vector_Alpha")
$ws7.Range("A2").Value = 0.5013171590436929
$ws7.Range("A3").Value = 0.09958128234312402
